$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 = 7)
$ws.Range("B2").Value = 0.1412068655553568
$ws.Range("C2").Value = 0.03888848024592173
$ws.Range("D2").Value = 0.2744609427486122
$ws.Range("E2").Value = -0.235231901775319
$ws.Range("F2").Value = 29.96816483908833
$ws.Range("G2").Value = 98.80967168531303
$ws.Range("H2").Value = 96.05520876326425

# Row 3 (A3 = 8)
$ws.Range("B3").Value = -0.8629830126873658
$ws.Range("C3").Value = 0.7619976890321466
$ws.Range("D3").Value = 1.174308772804542
$ws.Range("E3").Value = -0.803499379545246
$ws.Range("F3").Value = 69.60649528290904
$ws.Range("G3").Value = 97.96646863542702
$ws.Range("H3").Value = 95.12136927754187

# Row 4 (A4 = 9)
$ws.Range("B4").Value = 0.1466943608875539
$ws.Range("C4").Value = 0.2959403348166046
$ws.Range("D4").Value = -0.2303319100746397
$ws.Range("E4").Value = 0.6399546350019822
$ws.Range("F4").Value = 20.11364921925181
$ws.Range("G4").Value = 96.86018438497923
$ws.Range("H4").Value = 94.8521361022333

# Row 5 (A5 = 10)
$ws.Range("B5").Value = -0.5195417871016658
$ws.Range("C5").Value = 0.001403279013875681
$ws.Range("D5").Value = -0.1203563719034932
$ws.Range("E5").Value = -0.1552852937165759
$ws.Range("F5").Value = 3.070044701438284
$ws.Range("G5").Value = 98.60089019564199
$ws.Range("H5").Value = 97.08071733675607

# Row 6 (A6 = 11)
$ws.Range("B6").Value = 1.821489679462624
$ws.Range("C6").Value = -0.2716586800529331
$ws.Range("D6").Value = 0.0170193989029965
$ws.Range("E6").Value = -0.1441517209098913
$ws.Range("F6").Value = 3.748705388141238
$ws.Range("G6").Value = 98.15428500623874
$ws.Range("H6").Value = 94.78077109621393
